$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content corrections ("Correción Sprint y BD") ---
# Row 10 / Row 12: swap the "Responsable" names back to their correct owners
$ws.Range("D10").Value = "Fernando Abitia"
$ws.Range("D12").Value = "Joel Reyes"

# Row 11 / Row 12: rename tasks to reflect the modification work instead of the old CRUD wording
$ws.Range("B11").Value = "Modificar Jefe Departamento (Jefe)"
$ws.Range("B12").Value = "Modificar Tecnico (Técnico)"

# --- Restore the active selection left by the author when they saved ---
$ws.Range("H11:I11").Select()
